# Update "想去人数" (F column) counts on both the "展览" and "全部类型"
# worksheets, which contain duplicated data.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F.
$updates = @{
    3  = 701
    8  = 1622
    9  = 5753
    12 = 260
    13 = 74
    14 = 351
    15 = 123
    16 = 4684
    17 = 244
    18 = 1243
    20 = 96
    23 = 236
    24 = 84
    26 = 87
    27 = 369
    28 = 49
    30 = 74
    34 = 52
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
